$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.314.86"
$ws.Range("E2").Value = "  -2.70%  "
$ws.Range("D3").Value = "3.530.54"
$ws.Range("E3").Value = "  -4.54%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.24"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.85"
$ws.Range("E6").Value = "  -3.78%  "
$ws.Range("D7").Value = "3.525.37"
$ws.Range("E7").Value = "  -4.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.611"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -5.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.68"
$ws.Range("E11").Value = "  -2.95%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.587"
$ws.Range("E12").Value = "  -4.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.59"
$ws.Range("E13").Value = "  -3.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000275"
$ws.Range("E14").Value = "  -5.00%  "
$ws.Range("D15").Value = "4.092.30"
$ws.Range("E15").Value = "  -4.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.58"
$ws.Range("E16").Value = "  -5.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "630.20"
$ws.Range("E17").Value = "  -7.36%  "
$ws.Range("D18").Value = "3.534.96"
$ws.Range("E18").Value = "  -4.43%  "
$ws.Range("D19").Value = "69.279.01"
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.53"
$ws.Range("E21").Value = "  -2.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.25"
$ws.Range("E22").Value = "  -3.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.889"
$ws.Range("E23").Value = "  -5.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.03"
$ws.Range("E24").Value = "  -8.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.96"
$ws.Range("E25").Value = "  -4.28%  "
$ws.Range("E26").Value = "  -4.39%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.65"
$ws.Range("E28").Value = "  -7.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.31"
$ws.Range("E29").Value = "  -9.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.90"
$ws.Range("E30").Value = "  -6.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.16"
$ws.Range("E31").Value = "  -7.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.59"
$ws.Range("E32").Value = "  -6.22%  "
$ws.Range("E33").Value = "  -7.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.04"
$ws.Range("E34").Value = "  -6.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "632.58"
$ws.Range("E35").Value = "  +8.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.79"
$ws.Range("E36").Value = "  -3.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.53"
$ws.Range("E37").Value = "  -12.98%  "
$ws.Range("E38").Value = "  -4.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.38"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0455"
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("E42").Value = "  -5.61%  "
$ws.Range("D43").Value = "3.390.85"
$ws.Range("E43").Value = "  -6.52%  "
$ws.Range("E44").Value = "  -6.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "33.05"
$ws.Range("E45").Value = "  -7.58%  "
$ws.Range("E46").Value = "  -9.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.58"
$ws.Range("E47").Value = "  -7.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("E48").Value = "  -5.24%  "
$ws.Range("E49").Value = "  -2.54%  "
$ws.Range("E50").Value = "  +14.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.48"
$ws.Range("E51").Value = "  -2.65%  "
